$wb = $excel.ActiveWorkbook

# ALC!row64
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3536.1277
$ws.Cells.Item(64, 9).Value = 3485.1853
$ws.Cells.Item(64, 10).Value = 3604.9
$ws.Cells.Item(64, 11).Value = 3485.1853
$ws.Cells.Item(64, 12).Value = 3604.9
$ws.Cells.Item(64, 13).Value = -3237.1853
$ws.Cells.Item(64, 14).Value = -4100.9

# ALC!row67
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 3536.1277
$ws.Cells.Item(67, 9).Value = 3485.1853
$ws.Cells.Item(67, 10).Value = 3604.9
$ws.Cells.Item(67, 11).Value = 3485.1853
$ws.Cells.Item(67, 12).Value = 3604.9
$ws.Cells.Item(67, 13).Value = -2627.1853
$ws.Cells.Item(67, 14).Value = -5320.9

# ALC!row70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 8266.333000000001
$ws.Cells.Item(70, 9).Value = 9319.6
$ws.Cells.Item(70, 10).Value = 3000
$ws.Cells.Item(70, 11).Value = 27958.8
$ws.Cells.Item(70, 12).Value = 9000
$ws.Cells.Item(70, 13).Value = -27688.8
$ws.Cells.Item(70, 14).Value = -9540

# ALC!row73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 8266.333000000001
$ws.Cells.Item(73, 9).Value = 9319.6
$ws.Cells.Item(73, 10).Value = 3000
$ws.Cells.Item(73, 11).Value = 27958.8
$ws.Cells.Item(73, 12).Value = 9000
$ws.Cells.Item(73, 13).Value = -27022.8
$ws.Cells.Item(73, 14).Value = -10872

# ALC!row76
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3549.7
$ws.Cells.Item(76, 9).Value = 2501
$ws.Cells.Item(76, 11).Value = 2501
$ws.Cells.Item(76, 13).Value = -2186

# ALC!row79
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 3549.7
$ws.Cells.Item(79, 9).Value = 2501
$ws.Cells.Item(79, 11).Value = 2501
$ws.Cells.Item(79, 13).Value = -1409

# ALC!row106
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 4119.7334
$ws.Cells.Item(106, 9).Value = 2619.0908
$ws.Cells.Item(106, 11).Value = 2619.0908
$ws.Cells.Item(106, 13).Value = -1988.0908

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 46021.74
$ws.Cells.Item(137, 9).Value = 1883.7693
$ws.Cells.Item(137, 10).Value = 103401.1
$ws.Cells.Item(137, 11).Value = 5651.3079
$ws.Cells.Item(137, 12).Value = 310203.3
$ws.Cells.Item(137, 13).Value = -3101.3079
$ws.Cells.Item(137, 14).Value = -315303.3

# ARM!row44
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 14).ClearContents()

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1595.4333
$ws.Cells.Item(61, 9).Value = 1082.6538
$ws.Cells.Item(61, 10).Value = 4928.5
$ws.Cells.Item(61, 11).Value = 1082.6538
$ws.Cells.Item(61, 12).Value = 4928.5
$ws.Cells.Item(61, 13).Value = -870.6538
$ws.Cells.Item(61, 14).Value = -5352.5

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1828.6765
$ws.Cells.Item(132, 9).Value = 1754.2609
$ws.Cells.Item(132, 10).Value = 1984.2727
$ws.Cells.Item(132, 11).Value = 5262.7827
$ws.Cells.Item(132, 12).Value = 5952.8181
$ws.Cells.Item(132, 13).Value = -2732.7827
$ws.Cells.Item(132, 14).Value = -11012.8181

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1595.4333
$ws.Cells.Item(136, 9).Value = 1082.6538
$ws.Cells.Item(136, 10).Value = 4928.5
$ws.Cells.Item(136, 11).Value = 3247.9614
$ws.Cells.Item(136, 12).Value = 14785.5
$ws.Cells.Item(136, 13).Value = -697.9614000000001
$ws.Cells.Item(136, 14).Value = -19885.5

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1858.1489
$ws.Cells.Item(134, 9).Value = 1665.5428
$ws.Cells.Item(134, 10).Value = 2419.9167
$ws.Cells.Item(134, 11).Value = 4996.6284
$ws.Cells.Item(134, 12).Value = 7259.750100000001
$ws.Cells.Item(134, 13).Value = -2461.6284
$ws.Cells.Item(134, 14).Value = -12329.7501

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3451.4707
$ws.Cells.Item(31, 9).Value = 2014
$ws.Cells.Item(31, 10).Value = 5272.2666
$ws.Cells.Item(31, 11).Value = 2014
$ws.Cells.Item(31, 12).Value = 5272.2666
$ws.Cells.Item(31, 13).Value = -1719
$ws.Cells.Item(31, 14).Value = -5862.2666

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3451.4707
$ws.Cells.Item(34, 9).Value = 2014
$ws.Cells.Item(34, 10).Value = 5272.2666
$ws.Cells.Item(34, 11).Value = 2014
$ws.Cells.Item(34, 12).Value = 5272.2666
$ws.Cells.Item(34, 13).Value = -1812
$ws.Cells.Item(34, 14).Value = -5676.2666

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2613.25
$ws.Cells.Item(62, 9).Value = 2999.75
$ws.Cells.Item(62, 10).Value = 2226.75
$ws.Cells.Item(62, 11).Value = 2999.75
$ws.Cells.Item(62, 12).Value = 2226.75
$ws.Cells.Item(62, 13).Value = -2375.75
$ws.Cells.Item(62, 14).Value = -3474.75

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 2613.25
$ws.Cells.Item(65, 9).Value = 2999.75
$ws.Cells.Item(65, 10).Value = 2226.75
$ws.Cells.Item(65, 11).Value = 14998.75
$ws.Cells.Item(65, 12).Value = 11133.75
$ws.Cells.Item(65, 13).Value = -11878.75
$ws.Cells.Item(65, 14).Value = -17373.75

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1379766.9
$ws.Cells.Item(113, 9).Value = 6896980
$ws.Cells.Item(113, 10).Value = 463.65
$ws.Cells.Item(113, 11).Value = 20690940
$ws.Cells.Item(113, 12).Value = 1390.95
$ws.Cells.Item(113, 13).Value = -20688770
$ws.Cells.Item(113, 14).Value = -5730.95

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 4414.143
$ws.Cells.Item(131, 10).Value = 6188.9473
$ws.Cells.Item(131, 12).Value = 18566.8419
$ws.Cells.Item(131, 14).Value = -28646.8419

# CUL!row140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 4365196.5
$ws.Cells.Item(140, 9).Value = 6689834.5
$ws.Cells.Item(140, 11).Value = 20069503.5
$ws.Cells.Item(140, 13).Value = -20064323.5

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 9863.923000000001
$ws.Cells.Item(113, 9).Value = 875.7143
$ws.Cells.Item(113, 10).Value = 20350.166
$ws.Cells.Item(113, 11).Value = 875.7143
$ws.Cells.Item(113, 12).Value = 20350.166
$ws.Cells.Item(113, 13).Value = 1294.2857
$ws.Cells.Item(113, 14).Value = -24690.166

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4255.175
$ws.Cells.Item(132, 9).Value = 4333.75
$ws.Cells.Item(132, 10).Value = 3940.875
$ws.Cells.Item(132, 11).Value = 13001.25
$ws.Cells.Item(132, 12).Value = 11822.625
$ws.Cells.Item(132, 13).Value = -10471.25
$ws.Cells.Item(132, 14).Value = -16882.625

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2787.923
$ws.Cells.Item(68, 9).Value = 2774.5557
$ws.Cells.Item(68, 10).Value = 2818
$ws.Cells.Item(68, 11).Value = 2774.5557
$ws.Cells.Item(68, 12).Value = 2818
$ws.Cells.Item(68, 13).Value = -2025.5557
$ws.Cells.Item(68, 14).Value = -4316

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 2787.923
$ws.Cells.Item(71, 9).Value = 2774.5557
$ws.Cells.Item(71, 10).Value = 2818
$ws.Cells.Item(71, 11).Value = 13872.7785
$ws.Cells.Item(71, 12).Value = 14090
$ws.Cells.Item(71, 13).Value = -10128.7785
$ws.Cells.Item(71, 14).Value = -21578

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1614.09
$ws.Cells.Item(132, 9).Value = 1731.3678
$ws.Cells.Item(132, 10).Value = 829.2308
$ws.Cells.Item(132, 11).Value = 5194.1034
$ws.Cells.Item(132, 12).Value = 2487.6924
$ws.Cells.Item(132, 13).Value = -2664.1034
$ws.Cells.Item(132, 14).Value = -7547.6924

# WVR!row54
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 14).ClearContents()

# WVR!row113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 835.4
$ws.Cells.Item(113, 9).Value = 967.3333
$ws.Cells.Item(113, 10).Value = 637.5
$ws.Cells.Item(113, 11).Value = 2901.9999
$ws.Cells.Item(113, 12).Value = 1912.5
$ws.Cells.Item(113, 13).Value = -731.9998999999998
$ws.Cells.Item(113, 14).Value = -6252.5
